$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day 5 column header (F4): give the placeholder "Day 5" header its date ---
$ws.Range("F4").Value = "第五天`n日期:2025-11-26"

# --- Fill in Day 5 (column F) sleep-diary answers, previously blank ---
$ws.Range("F5").Value  = "7：50"   # 您今天早上几点醒来?
$ws.Range("F6").Value  = "7：50"   # 您今天几点起床?
$ws.Range("F7").Value  = "22：40"  # 您昨晚几点上床?

# --- Correct the "24:00" notation to "00:00" for day 2 and day 3 lights-off times ---
$ws.Range("C8").Value = "00：00"
$ws.Range("D8").Value = "00：00"

$ws.Range("F8").Value  = "23：20"  # 您昨晚几点熄灯?
$ws.Range("F9").Value  = 50        # 入睡用时(分钟)
$ws.Range("F10").Value = 3         # 整晚醒来几次
$ws.Range("F11").Value = 20        # 整晚醒了多长时间(分钟)
$ws.Range("F12").Value = 420       # 整晚睡了多长时间(分钟)
$ws.Range("F13").Value = "无"      # 是否使用影响睡眠的物质
$ws.Range("F14").Value = 90        # 睡前使用电子产品时长(分钟)
$ws.Range("F15").Value = 3         # 睡眠质量
$ws.Range("F16").Value = 1         # 睡前身体紧张程度
$ws.Range("F17").Value = 4         # 睡前精神紧张程度
$ws.Range("F18").Value = "无"      # 白天是否小睡

# --- Update the view state to reflect the user scrolling to / selecting the newly edited cell ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F18").Select()
